$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 6993121
$ws.Range("I8").Value = 6993121
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 20979363
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -20979224
$ws.Range("N8").ClearContents()

$ws.Range("H39").Value = 531.92
$ws.Range("I39").Value = 76.27273
$ws.Range("J39").Value = 889.9286
$ws.Range("K39").Value = 228.81819
$ws.Range("L39").Value = 2669.7858
$ws.Range("M39").Value = 67.18181000000001
$ws.Range("N39").Value = -3261.7858

$ws.Range("H105").Value = 29400
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 29400
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 29400
$ws.Range("N105").Value = -36388

$ws.Range("H113").Value = 3558.2273
$ws.Range("I113").Value = 2609.5833
$ws.Range("J113").Value = 4696.6
$ws.Range("K113").Value = 2609.5833
$ws.Range("L113").Value = 4696.6
$ws.Range("M113").Value = 644.4167000000002
$ws.Range("N113").Value = -11204.6

$ws.Range("H137").Value = 1810
$ws.Range("I137").Value = 3388.1667
$ws.Range("J137").Value = 986.6087
$ws.Range("K137").Value = 10164.5001
$ws.Range("L137").Value = 2959.8261
$ws.Range("M137").Value = -7614.500100000001
$ws.Range("N137").Value = -8059.8261

$ws.Range("H138").Value = 1827.6962
$ws.Range("I138").Value = 1361.375
$ws.Range("J138").Value = 2305.9744
$ws.Range("K138").Value = 4084.125
$ws.Range("L138").Value = 6917.9232
$ws.Range("M138").Value = 1055.875
$ws.Range("N138").Value = -17197.9232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6861.82
$ws.Range("I32").Value = 5426.449
$ws.Range("J32").Value = 18475.273
$ws.Range("K32").Value = 5426.449
$ws.Range("L32").Value = 18475.273
$ws.Range("M32").Value = -5139.449
$ws.Range("N32").Value = -19049.273

$ws.Range("H101").Value = 29500
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 29500
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 29500
$ws.Range("N101").Value = -35990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 30271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30271
$ws.Range("N63").Value = -31643

$ws.Range("H66").Value = 30271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90813
$ws.Range("N66").Value = -97677

$ws.Range("H103").Value = 12776.167
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 12776.167
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 12776.167
$ws.Range("N103").Value = -15120.167

$ws.Range("H134").Value = 4387.8
$ws.Range("I134").Value = 2009.3043
$ws.Range("J134").Value = 6097.3438
$ws.Range("K134").Value = 6027.9129
$ws.Range("L134").Value = 18292.0314
$ws.Range("M134").Value = -3492.9129
$ws.Range("N134").Value = -23362.0314

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3077.2354
$ws.Range("I16").Value = 2524.077
$ws.Range("J16").Value = 4875
$ws.Range("K16").Value = 2524.077
$ws.Range("L16").Value = 4875
$ws.Range("M16").Value = -2237.077
$ws.Range("N16").Value = -5449

$ws.Range("H86").Value = 3793.7585
$ws.Range("I86").Value = 3039.9565
$ws.Range("J86").Value = 6683.3335
$ws.Range("K86").Value = 3039.9565
$ws.Range("L86").Value = 6683.3335
$ws.Range("M86").Value = -1916.9565
$ws.Range("N86").Value = -8929.333500000001

$ws.Range("H89").Value = 3793.7585
$ws.Range("I89").Value = 3039.9565
$ws.Range("J89").Value = 6683.3335
$ws.Range("K89").Value = 15199.7825
$ws.Range("L89").Value = 33416.6675
$ws.Range("M89").Value = -9583.782499999999
$ws.Range("N89").Value = -44648.6675

$ws.Range("H93").Value = 9380
$ws.Range("I93").Value = 5370
$ws.Range("J93").Value = 17400
$ws.Range("K93").Value = 5370
$ws.Range("L93").Value = 17400
$ws.Range("M93").Value = -3498
$ws.Range("N93").Value = -21144

$ws.Range("H105").Value = 1420.8462
$ws.Range("I105").Value = 1496.3636
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1496.3636
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 250.6364000000001
$ws.Range("N105").Value = -4499.5

$ws.Range("H107").Value = 878.13043
$ws.Range("I107").Value = 312.46667
$ws.Range("J107").Value = 1938.75
$ws.Range("K107").Value = 312.46667
$ws.Range("L107").Value = 1938.75
$ws.Range("M107").Value = 1607.53333
$ws.Range("N107").Value = -5778.75

$ws.Range("H113").Value = 3077.2354
$ws.Range("I113").Value = 2524.077
$ws.Range("J113").Value = 4875
$ws.Range("K113").Value = 2524.077
$ws.Range("L113").Value = 4875
$ws.Range("M113").Value = -354.0770000000002
$ws.Range("N113").Value = -9215

$ws.Range("H122").Value = 2389.25
$ws.Range("I122").Value = 775
$ws.Range("J122").Value = 4003.5
$ws.Range("K122").Value = 2325
$ws.Range("L122").Value = 12010.5
$ws.Range("M122").Value = 125

$ws.Range("H134").Value = 2218.9375
$ws.Range("I134").Value = 1353.8334
$ws.Range("J134").Value = 2738
$ws.Range("K134").Value = 4061.5002
$ws.Range("L134").Value = 8214
$ws.Range("M134").Value = -1526.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 58.583332
$ws.Range("I14").Value = 58.583332
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 175.749996
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -2.74999600000001

$ws.Range("H75").Value = 9000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 9000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 27000
$ws.Range("N75").Value = -28996
$ws.Range("M75").ClearContents()

$ws.Range("H78").Value = 9000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 9000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 81000
$ws.Range("N78").Value = -90984
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2139.0625
$ws.Range("I16").Value = 1392.3684
$ws.Range("J16").Value = 3230.3845
$ws.Range("K16").Value = 1392.3684
$ws.Range("L16").Value = 3230.3845
$ws.Range("M16").Value = -1222.3684
$ws.Range("N16").Value = -3570.3845

$ws.Range("H40").Value = 4559.8887
$ws.Range("I40").Value = 3451.2
$ws.Range("J40").Value = 5945.75
$ws.Range("K40").Value = 3451.2
$ws.Range("L40").Value = 5945.75
$ws.Range("M40").Value = -3315.2
$ws.Range("N40").Value = -6217.75

$ws.Range("H97").Value = 12629.071
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 12629.071
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 12629.071
$ws.Range("N97").Value = -14611.071

$ws.Range("H103").Value = 38643.285
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 38643.285
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 38643.285
$ws.Range("N103").Value = -40987.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 18293.084
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 18293.084
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 18293.084
$ws.Range("N97").Value = -20275.084

$ws.Range("H107").Value = 758.3125
$ws.Range("I107").Value = 692.75
$ws.Range("J107").Value = 955
$ws.Range("K107").Value = 2078.25
$ws.Range("L107").Value = 2865
$ws.Range("M107").Value = -158.25
$ws.Range("N107").Value = -6705
